# Applies the cryptos.xlsx price/volume refresh described in the commit:
# "Updated cryptos list on Fri Jul  5 11:23:33 UTC 2024 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "55.277.45"
$ws.Range("E2").Value = "  -4.38%  "

# Row 3
$ws.Range("D3").Value = "2.935.05"
$ws.Range("E3").Value = "  -7.23%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.10%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "479.50"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -8.98%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "129.10"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.47%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.15%  "

# Row 8
$ws.Range("D8").Value = "2.929.50"
$ws.Range("E8").Value = "  -7.38%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.413"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -8.90%  "

# Row 10
$ws.Range("E10").Value = "  -6.67%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0993"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -10.65%  "

# Row 12
$ws.Range("E12").Value = "  -12.17%  "

# Row 13
$ws.Range("E13").Value = "  -1.28%  "

# Row 14
$ws.Range("D14").Value = "3.437.91"
$ws.Range("E14").Value = "  -7.30%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "24.10"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -6.95%  "

# Row 16
$ws.Range("D16").Value = "55.211.58"
$ws.Range("E16").Value = "  -4.46%  "

# Row 17
$ws.Range("D17").Value = "2.935.16"
$ws.Range("E17").Value = "  -7.45%  "

# Row 18
$ws.Range("E18").Value = "  -10.20%  "

# Row 19
$ws.Range("E19").Value = "  -5.52%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.78"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -10.08%  "

# Row 21
$ws.Range("E21").Value = "  -9.36%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "307.51"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -11.40%  "

# Row 23
$ws.Range("E23").Value = "  +0.02%  "

# Row 24
$ws.Range("E24").Value = "  -11.53%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "59.41"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -14.76%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.998"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.14%  "

# Row 27
$ws.Range("E27").Value = "  -6.86%  "

# Row 28
$ws.Range("E28").Value = "  +0.08%  "

# Row 29
$ws.Range("D29").Value = "0.0₃0832"
$ws.Range("E29").Value = "  -13.59%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.41"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -7.00%  "

# Row 31
$ws.Range("B31").Value = "Fetch.AI"
$ws.Range("C31").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.15"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -5.38%  "

# Row 32
$ws.Range("B32").Value = "InternetComputer(DFINITY)"
$ws.Range("C32").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.39"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -7.80%  "

# Row 33
$ws.Range("E33").Value = "  -11.93%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "19.03"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -12.61%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "146.81"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -8.40%  "

# Row 36
$ws.Range("E36").Value = "  -13.08%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.55"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -11.37%  "

# Row 38
$ws.Range("E38").Value = "  -11.46%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "23.27"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -10.38%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0635"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -9.10%  "

# Row 41
$ws.Range("D41").Value = "2.965.86"
$ws.Range("E41").Value = "  -7.08%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.00"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.11%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "35.72"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -12.13%  "

# Row 44
$ws.Range("E44").Value = "  -9.18%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.623"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -11.12%  "

# Row 46
$ws.Range("E46").Value = "  -7.56%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.48"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -12.30%  "

# Row 48
$ws.Range("D48").Value = "2.098.11"
$ws.Range("E48").Value = "  -7.77%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0226"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -4.70%  "

# Row 50
$ws.Range("B50").Value = "InjectiveProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "18.58"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -9.72%  "

# Row 51
$ws.Range("B51").Value = "Cosmos"
$ws.Range("C51").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "5.48"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -11.46%  "
